# Commit: Thu, Jul 02, 2020  3:06:08 PM
#
# 1) Slide 6 ("SOURCES OF FINANCE") table: switch the table's style from
#    the custom "Table_0" style to the built-in "Medium Style 2 - Accent 1"
#    table style (GUID {337D9D03-21C2-4685-AC51-E491D688D4C3}).
#
# 2) The deck's theme colour scheme (currently the "Integral" palette) is
#    changed over to the stock "Office Theme" palette - dk2/lt2/accent1-6/
#    hlink/folHlink all move to the default Office RGB values (dk1/lt1
#    stay black/white, and the font + format schemes are already shared
#    between the two themes, so only the colours need to move).

$p = $ppt.ActivePresentation

# --- 1) table style on slide 6 -------------------------------------------
$s6 = $p.Slides.Item(6)
for ($i = 1; $i -le $s6.Shapes.Count; $i++) {
    $shp = $s6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{337D9D03-21C2-4685-AC51-E491D688D4C3}")
    }
}

# --- 2) theme colours: Integral -> Office Theme ---------------------------
$themeColors = $p.Slides.Item(1).ThemeColorScheme

# index: 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3 8=accent4
#        9=accent5 10=accent6 11=hlink 12=folHlink
# dk1/lt1 (black/white) are identical between the two themes, so only
# indices 3-12 need updating.
$officeRgb = @{
    3  = 6968388    # dk2      44546A
    4  = 15132391   # lt2      E7E6E6
    5  = 13998939   # accent1  5B9BD5
    6  = 3243501     # accent2  ED7D31
    7  = 10855845    # accent3  A5A5A5
    8  = 49407       # accent4  FFC000
    9  = 12874308    # accent5  4472C4
    10 = 4697456      # accent6  70AD47
    11 = 12673797     # hlink    0563C1
    12 = 7491477      # folHlink 954F72
}

foreach ($idx in $officeRgb.Keys) {
    $themeColors.Item($idx).RGB = $officeRgb[$idx]
}
